$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Fix the typo in the shared string "cClasTrib" -> "cClassTrib"
#    (this cell currently lives at D6, before the row-insert below).
# ------------------------------------------------------------------
$ws.Range("D6").Value = "cClassTrib"

# ------------------------------------------------------------------
# 2) Turn off gridlines for the sheet view.
# ------------------------------------------------------------------
$excel.ActiveWindow.DisplayGridlines = $false

# ------------------------------------------------------------------
# 3) Insert two new blank rows right before the header row (row 6),
#    pushing the header from row 6 down to row 8.
# ------------------------------------------------------------------
$ws.Rows.Item(6).Resize(2).Insert()

# ------------------------------------------------------------------
# 4) Stamp the blank "logo" column cells (A1:A5) and B3 so they exist
#    explicitly in the sheet (part of the reworked layout).
# ------------------------------------------------------------------
$ws.Range("A1:A5").HorizontalAlignment = -4131
$ws.Range("A1:A5").HorizontalAlignment = 1
$ws.Range("B3").HorizontalAlignment = -4131
$ws.Range("B3").HorizontalAlignment = 1

# ------------------------------------------------------------------
# 5) New blank spacer rows (6 and 7): right-align + merge each row.
# ------------------------------------------------------------------
$spacer = $ws.Range("A6:F7")
$spacer.HorizontalAlignment = -4152
$ws.Range("A6:F6").Merge()
$ws.Range("A7:F7").Merge()

# ------------------------------------------------------------------
# 6) Explicit row heights (forces customHeight="true" on the rows
#    that need it; rows 4 and 5 are intentionally left on auto).
# ------------------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 24
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(3).RowHeight = 15
$ws.Rows.Item(6).RowHeight = 15
$ws.Rows.Item(7).RowHeight = 15
$ws.Rows.Item(8).RowHeight = 15

# ------------------------------------------------------------------
# 7) Move the active selection to D11 (was D15).
# ------------------------------------------------------------------
$ws.Range("D11").Select()

# ------------------------------------------------------------------
# 8) Picture fixes: rename, drop the near-zero rotation and correct
#    the position / size to the exact target geometry.
# ------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Name = "Picture 1"
$shp.Rotation = 0
$shp.Left = 0.08503937007874016
$shp.Top = 0.11338582677165354
$shp.Width = 140.54173228346457
$shp.Height = 91.75748031496063

Write-Output "edit complete"
